$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2508241.8
$ws.Range("I132").Value = 2917369
$ws.Range("J132").Value = 2337.75
$ws.Range("K132").Value = 8752107
$ws.Range("L132").Value = 7013.25
$ws.Range("M132").Value = -8749577
$ws.Range("N132").Value = -12073.25

$ws.Range("H135").Value = 3435.8298
$ws.Range("I135").Value = 527.5833
$ws.Range("J135").Value = 12953.728
$ws.Range("K135").Value = 4748.2497
$ws.Range("L135").Value = 116583.552
$ws.Range("M135").Value = -2213.2497
$ws.Range("N135").Value = -121653.552

$ws.Range("H137").Value = 41668430
$ws.Range("I137").Value = 1365.2106
$ws.Range("K137").Value = 4095.6318
$ws.Range("M137").Value = -1545.6318

$ws.Range("H138").Value = 3250.0156
$ws.Range("I138").Value = 2704.7896
$ws.Range("J138").Value = 3480.2222
$ws.Range("K138").Value = 8114.3688
$ws.Range("L138").Value = 10440.6666
$ws.Range("M138").Value = -2974.3688
$ws.Range("N138").Value = -20720.6666

$ws.Range("H141").Value = 1977.1774
$ws.Range("I141").Value = 1342.375
$ws.Range("J141").Value = 3131.3635
$ws.Range("K141").Value = 4027.125
$ws.Range("L141").Value = 9394.0905
$ws.Range("M141").Value = 1152.875
$ws.Range("N141").Value = -19754.0905

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 43479348
$ws.Range("I2").Value = 142857680
$ws.Range("J2").Value = 1331.0625
$ws.Range("K2").Value = 142857680
$ws.Range("L2").Value = 1331.0625
$ws.Range("M2").Value = -142857567
$ws.Range("N2").Value = -1557.0625

$ws.Range("H32").Value = 17453.72
$ws.Range("I32").Value = 17537.795
$ws.Range("J32").Value = 16955.285
$ws.Range("K32").Value = 17537.795
$ws.Range("L32").Value = 16955.285
$ws.Range("M32").Value = -17250.795
$ws.Range("N32").Value = -17529.285

$ws.Range("H45").Value = 41667756
$ws.Range("I45").Value = 83334230
$ws.Range("J45").Value = 1275
$ws.Range("K45").Value = 83334230
$ws.Range("L45").Value = 1275
$ws.Range("M45").Value = -83333853
$ws.Range("N45").Value = -2029

$ws.Range("H116").Value = 43479348
$ws.Range("I116").Value = 142857680
$ws.Range("J116").Value = 1331.0625
$ws.Range("K116").Value = 142857680
$ws.Range("L116").Value = 1331.0625
$ws.Range("M116").Value = -142855386
$ws.Range("N116").Value = -5919.0625

$ws.Range("H122").Value = 2840.25
$ws.Range("I122").Value = 2840.25
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8520.75
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -6070.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 43479348
$ws.Range("I3").Value = 142857680
$ws.Range("J3").Value = 1331.0625
$ws.Range("K3").Value = 142857680
$ws.Range("L3").Value = 1331.0625
$ws.Range("M3").Value = -142857566
$ws.Range("N3").Value = -1559.0625

$ws.Range("H55").Value = 59780
$ws.Range("J55").Value = 59780
$ws.Range("L55").Value = 59780
$ws.Range("N55").Value = -60326

$ws.Range("H134").Value = 5556.4
$ws.Range("I134").Value = 7366.5884
$ws.Range("J134").Value = 3189.2307
$ws.Range("K134").Value = 22099.7652
$ws.Range("L134").Value = 9567.6921
$ws.Range("M134").Value = -19564.7652
$ws.Range("N134").Value = -14637.6921

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1094.3125
$ws.Range("I58").Value = 1128.0889
$ws.Range("J58").Value = 587.6667
$ws.Range("K58").Value = 1128.0889
$ws.Range("L58").Value = 587.6667
$ws.Range("M58").Value = -925.0889
$ws.Range("N58").Value = -993.6667

$ws.Range("H132").Value = 4631138
$ws.Range("I132").Value = 1227.6364
$ws.Range("J132").Value = 25002742
$ws.Range("K132").Value = 3682.9092
$ws.Range("L132").Value = 75008226
$ws.Range("M132").Value = -1152.9092
$ws.Range("N132").Value = -75013286

$ws.Range("H136").Value = 1094.3125
$ws.Range("I136").Value = 1128.0889
$ws.Range("J136").Value = 587.6667
$ws.Range("K136").Value = 3384.2667
$ws.Range("L136").Value = 1763.0001
$ws.Range("M136").Value = -834.2667000000001
$ws.Range("N136").Value = -6863.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 11765720
$ws.Range("I26").Value = 60
$ws.Range("J26").Value = 15385923
$ws.Range("K26").Value = 180
$ws.Range("L26").Value = 46157769
$ws.Range("M26").Value = 108
$ws.Range("N26").Value = -46158345

$ws.Range("H80").Value = 3846
$ws.Range("J80").Value = 4320.3335
$ws.Range("L80").Value = 12961.0005
$ws.Range("N80").Value = -14833.0005

$ws.Range("H83").Value = 3846
$ws.Range("J83").Value = 4320.3335
$ws.Range("L83").Value = 38883.0015
$ws.Range("N83").Value = -48243.0015

$ws.Range("H131").Value = 2060198.9
$ws.Range("J131").Value = 2469939.5
$ws.Range("L131").Value = 7409818.5
$ws.Range("N131").Value = -7419898.5

$ws.Range("H132").Value = 1982.0435
$ws.Range("I132").Value = 1053
$ws.Range("J132").Value = 2177.6316
$ws.Range("K132").Value = 9477
$ws.Range("L132").Value = 19598.6844
$ws.Range("M132").Value = -6947
$ws.Range("N132").Value = -24658.6844

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4291.718
$ws.Range("I132").Value = 4543.0884
$ws.Range("J132").Value = 2582.4
$ws.Range("K132").Value = 13629.2652
$ws.Range("L132").Value = 7747.200000000001
$ws.Range("M132").Value = -11099.2652
$ws.Range("N132").Value = -12807.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 45456948
$ws.Range("I7").Value = 1900
$ws.Range("J7").Value = 62502588
$ws.Range("K7").Value = 1900
$ws.Range("L7").Value = 62502588
$ws.Range("M7").Value = -1788
$ws.Range("N7").Value = -62502812

$ws.Range("H126").Value = 45456948
$ws.Range("I126").Value = 1900
$ws.Range("J126").Value = 62502588
$ws.Range("K126").Value = 5700
$ws.Range("L126").Value = 187507764
$ws.Range("M126").Value = -3230
$ws.Range("N126").Value = -187512704

$ws.Range("H127").Value = 60000
$ws.Range("J127").Value = 60000
$ws.Range("L127").Value = 60000
$ws.Range("N127").Value = -69920

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 13889271
$ws.Range("I107").Value = 5682198
$ws.Range("J107").Value = 30303416
$ws.Range("K107").Value = 17046594
$ws.Range("L107").Value = 90910248
$ws.Range("M107").Value = -17044674
$ws.Range("N107").Value = -90914088

$ws.Range("H136").Value = 2651.037
$ws.Range("I136").Value = 2698.7778
$ws.Range("J136").Value = 2555.5557
$ws.Range("K136").Value = 8096.3334
$ws.Range("L136").Value = 7666.6671
$ws.Range("M136").Value = -5546.3334
$ws.Range("N136").Value = -12766.6671
